# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh handback report
# generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-06 15:28:21"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-06 15:27:56"
$wsZhCn.Range("K2").Value = "2016-09-06 15:29:15"

# de-de sheet: Correspond Handoff Datetime (shares text with Overview G2) /
# Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-06 15:28:21"
$wsDeDe.Range("K2").Value = "2016-09-06 15:29:36"

$wb.Save()
